$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: "Gestione" activity with a note and 1 hour logged.
$ws.Range("B9").Value = "Gestione"
$ws.Range("C9").Value = "Appunti e cose varie"
$ws.Range("D9").Value = 1/24

# Matches the author's final selection position after adding the row.
[void]$ws.Range("E13").Select()
